$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "time_taken" header in F1, matching the style of the other header cells (E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate F2:F11 with the time_taken values recorded for each panel row
$ws.Range("F2").Value = "2021-10-05 10:50:23.238517"
$ws.Range("F3").Value = "2021-10-05 10:50:23.238528"
$ws.Range("F4").Value = "2021-10-05 10:50:23.238532"
$ws.Range("F5").Value = "2021-10-05 10:50:23.238534"
$ws.Range("F6").Value = "2021-10-05 10:50:23.238537"
$ws.Range("F7").Value = "2021-10-05 10:50:23.238540"
$ws.Range("F8").Value = "2021-10-05 10:50:23.238543"
$ws.Range("F9").Value = "2021-10-05 10:50:23.238546"
$ws.Range("F10").Value = "2021-10-05 10:50:23.238548"
$ws.Range("F11").Value = "2021-10-05 10:50:23.238551"
